# Regenerate orders with updated distance/size codes.
#
# The experiment's distance codes change D64 -> D69, D51 -> D55, D80 -> D86,
# and the "large" size code changes S30 -> S31 (S20 and S25 stay the same).
# These codes appear embedded inside several text columns (Condition,
# Filename_Left, Filename_Right, Distance, Size), so every cell in those
# columns needs the substring(s) updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Columns containing text with embedded Distance/Size codes:
#   B = Condition, D = Filename_Left, E = Filename_Right,
#   H = Distance,  J = Size
$cols = @(2, 4, 5, 8, 10)

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $orig = $cell.Text

        $new = $orig -replace "D64", "D69"
        $new = $new -replace "D51", "D55"
        $new = $new -replace "D80", "D86"
        $new = $new -replace "S30", "S31"

        if ($new -ne $orig) {
            $cell.Value = $new
        }
    }
}
